# Update the dSF ("F") column values for Sheet1 (harvey_matt.xlsx)
# per the "repull data, push all data, mean calculation" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    3  = -2
    4  = -1
    5  = -1
    6  = -3
    7  = -3
    8  = 1
    9  = 4
    10 = 2
    11 = 4
    12 = -4
    13 = 2
    14 = -2
    15 = -2
    16 = -2
    17 = -6
    18 = -1
    19 = -2
    20 = -1
    21 = -6
    22 = -5
    23 = -4
    24 = 8
    25 = 1
    28 = 1
    29 = 2
    30 = -1
    33 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
